$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("E2").Value = 3
$ws.Range("G2").Value = 1.476790666666667
$ws.Range("H2").Value = 4.430372
$ws.Range("I2").Value = 0.01966406119907831
$ws.Range("J2").Value = 0.0196640611990783
$ws.Range("K2").Value = 3
$ws.Range("M2").Value = 23.63579766666667
$ws.Range("N2").Value = 70.907393
$ws.Range("O2").Value = 0.06827844587621175
$ws.Range("P2").Value = 0.06827844587621175
$ws.Range("Q2").Value = 34.90512539335511
$ws.Range("R2").Value = 314.146128540196
$ws.Range("S2").Value = 0.001342631538287784
$ws.Range("T2").Value = 0.001342631538287783
$ws.Range("E3").Value = 3
$ws.Range("G3").Value = 1.476790666666667
$ws.Range("H3").Value = 4.430372
$ws.Range("I3").Value = 0.01966406119907831
$ws.Range("J3").Value = 0.0196640611990783
$ws.Range("K3").Value = 3
$ws.Range("M3").Value = 181.2883913333334
$ws.Range("N3").Value = 543.865174
$ws.Range("O3").Value = 0.5237009467675041
$ws.Range("P3").Value = 0.523700946767504
$ws.Range("Q3").Value = 267.725004296081
$ws.Range("R3").Value = 2409.525038664728
$ws.Range("S3").Value = 0.01029808746725145
$ws.Range("T3").Value = 0.01029808746725145
$ws.Range("E4").Value = 3
$ws.Range("G4").Value = 1.476790666666667
$ws.Range("H4").Value = 4.430372
$ws.Range("I4").Value = 0.01966406119907831
$ws.Range("J4").Value = 0.0196640611990783
$ws.Range("K4").Value = 3
$ws.Range("M4").Value = 111.1005463333333
$ws.Range("N4").Value = 333.301639
$ws.Range("O4").Value = 0.3209442197221123
$ws.Range("P4").Value = 0.3209442197221123
$ws.Range("Q4").Value = 164.0722498866342
$ws.Range("R4").Value = 1476.650248979708
$ws.Range("S4").Value = 0.006311066778106051
$ws.Range("T4").Value = 0.00631106677810605
$ws.Range("E5").Value = 3
$ws.Range("G5").Value = 1.476790666666667
$ws.Range("H5").Value = 4.430372
$ws.Range("I5").Value = 0.01966406119907831
$ws.Range("J5").Value = 0.0196640611990783
$ws.Range("K5").Value = 3
$ws.Range("M5").Value = 30.14303933333333
$ws.Range("N5").Value = 90.42911799999999
$ws.Range("O5").Value = 0.08707638763417187
$ws.Range("P5").Value = 0.08707638763417187
$ws.Range("Q5").Value = 44.51495915243289
$ws.Range("R5").Value = 400.634632371896
$ws.Range("S5").Value = 0.001712275415433021
$ws.Range("T5").Value = 0.001712275415433021
$ws.Range("E6").Value = 3
$ws.Range("G6").Value = 44.80640933333333
$ws.Range("H6").Value = 134.419228
$ws.Range("I6").Value = 0.5966153464595884
$ws.Range("J6").Value = 0.5966153464595884
$ws.Range("K6").Value = 3
$ws.Range("M6").Value = 23.63579766666667
$ws.Range("N6").Value = 70.907393
$ws.Range("O6").Value = 0.06827844587621175
$ws.Range("P6").Value = 0.06827844587621175
$ws.Range("Q6").Value = 1059.035225172511
$ws.Range("R6").Value = 9531.317026552602
$ws.Range("S6").Value = 0.04073596864215832
$ws.Range("T6").Value = 0.04073596864215832
$ws.Range("E7").Value = 3
$ws.Range("G7").Value = 44.80640933333333
$ws.Range("H7").Value = 134.419228
$ws.Range("I7").Value = 0.5966153464595884
$ws.Range("J7").Value = 0.5966153464595884
$ws.Range("K7").Value = 3
$ws.Range("M7").Value = 181.2883913333334
$ws.Range("N7").Value = 543.865174
$ws.Range("O7").Value = 0.5237009467675041
$ws.Range("P7").Value = 0.523700946767504
$ws.Range("Q7").Value = 8122.881869462853
$ws.Range("R7").Value = 73105.93682516566
$ws.Range("S7").Value = 0.3124480217969089
$ws.Range("T7").Value = 0.3124480217969088
$ws.Range("E8").Value = 3
$ws.Range("G8").Value = 44.80640933333333
$ws.Range("H8").Value = 134.419228
$ws.Range("I8").Value = 0.5966153464595884
$ws.Range("J8").Value = 0.5966153464595884
$ws.Range("K8").Value = 3
$ws.Range("M8").Value = 111.1005463333333
$ws.Range("N8").Value = 333.301639
$ws.Range("O8").Value = 0.3209442197221123
$ws.Range("P8").Value = 0.3209442197221123
$ws.Range("Q8").Value = 4978.016556168298
$ws.Range("R8").Value = 44802.14900551468
$ws.Range("S8").Value = 0.1914802468437103
$ws.Range("T8").Value = 0.1914802468437103
$ws.Range("E9").Value = 3
$ws.Range("G9").Value = 44.80640933333333
$ws.Range("H9").Value = 134.419228
$ws.Range("I9").Value = 0.5966153464595884
$ws.Range("J9").Value = 0.5966153464595884
$ws.Range("K9").Value = 3
$ws.Range("M9").Value = 30.14303933333333
$ws.Range("N9").Value = 90.42911799999999
$ws.Range("O9").Value = 0.08707638763417187
$ws.Range("P9").Value = 0.08707638763417187
$ws.Range("Q9").Value = 1350.6013589201
$ws.Range("R9").Value = 12155.4122302809
$ws.Range("S9").Value = 0.05195110917681087
$ws.Range("T9").Value = 0.05195110917681087
$ws.Range("E10").Value = 3
$ws.Range("G10").Value = 2.566584666666667
$ws.Range("H10").Value = 7.699754
$ws.Range("I10").Value = 0.03417510626056863
$ws.Range("J10").Value = 0.03417510626056863
$ws.Range("K10").Value = 3
$ws.Range("M10").Value = 23.63579766666667
$ws.Range("N10").Value = 70.907393
$ws.Range("O10").Value = 0.06827844587621175
$ws.Range("P10").Value = 0.06827844587621175
$ws.Range("Q10").Value = 60.66327587570245
$ws.Range("R10").Value = 545.9694828813221
$ws.Range("S10").Value = 0.002333423143126021
$ws.Range("T10").Value = 0.002333423143126021
$ws.Range("E11").Value = 3
$ws.Range("G11").Value = 2.566584666666667
$ws.Range("H11").Value = 7.699754
$ws.Range("I11").Value = 0.03417510626056863
$ws.Range("J11").Value = 0.03417510626056863
$ws.Range("K11").Value = 3
$ws.Range("M11").Value = 181.2883913333334
$ws.Range("N11").Value = 543.865174
$ws.Range("O11").Value = 0.5237009467675041
$ws.Range("P11").Value = 0.523700946767504
$ws.Range("Q11").Value = 465.2920054407996
$ws.Range("R11").Value = 4187.628048967196
$ws.Range("S11").Value = 0.01789753550453985
$ws.Range("T11").Value = 0.01789753550453985
$ws.Range("E12").Value = 3
$ws.Range("G12").Value = 2.566584666666667
$ws.Range("H12").Value = 7.699754
$ws.Range("I12").Value = 0.03417510626056863
$ws.Range("J12").Value = 0.03417510626056863
$ws.Range("K12").Value = 3
$ws.Range("M12").Value = 111.1005463333333
$ws.Range("N12").Value = 333.301639
$ws.Range("O12").Value = 0.3209442197221123
$ws.Range("P12").Value = 0.3209442197221123
$ws.Range("Q12").Value = 285.1489586774229
$ws.Range("R12").Value = 2566.340628096806
$ws.Range("S12").Value = 0.01096830281271848
$ws.Range("T12").Value = 0.01096830281271848
$ws.Range("E13").Value = 3
$ws.Range("G13").Value = 2.566584666666667
$ws.Range("H13").Value = 7.699754
$ws.Range("I13").Value = 0.03417510626056863
$ws.Range("J13").Value = 0.03417510626056863
$ws.Range("K13").Value = 3
$ws.Range("M13").Value = 30.14303933333333
$ws.Range("N13").Value = 90.42911799999999
$ws.Range("O13").Value = 0.08707638763417187
$ws.Range("P13").Value = 0.08707638763417187
$ws.Range("Q13").Value = 77.36466255966356
$ws.Range("R13").Value = 696.281963036972
$ws.Range("S13").Value = 0.002975844800184288
$ws.Range("T13").Value = 0.002975844800184288
$ws.Range("E14").Value = 3
$ws.Range("G14").Value = 26.25121566666667
$ws.Range("H14").Value = 78.753647
$ws.Range("I14").Value = 0.3495454860807646
$ws.Range("J14").Value = 0.3495454860807646
$ws.Range("K14").Value = 3
$ws.Range("M14").Value = 23.63579766666667
$ws.Range("N14").Value = 70.907393
$ws.Range("O14").Value = 0.06827844587621175
$ws.Range("P14").Value = 0.06827844587621175
$ws.Range("Q14").Value = 620.4684220013634
$ws.Range("R14").Value = 5584.215798012271
$ws.Range("S14").Value = 0.02386642255263962
$ws.Range("T14").Value = 0.02386642255263962
$ws.Range("E15").Value = 3
$ws.Range("G15").Value = 26.25121566666667
$ws.Range("H15").Value = 78.753647
$ws.Range("I15").Value = 0.3495454860807646
$ws.Range("J15").Value = 0.3495454860807646
$ws.Range("K15").Value = 3
$ws.Range("M15").Value = 181.2883913333334
$ws.Range("N15").Value = 543.865174
$ws.Range("O15").Value = 0.5237009467675041
$ws.Range("P15").Value = 0.5237009467675041
$ws.Range("Q15").Value = 4759.040658754398
$ws.Range("R15").Value = 42831.36592878958
$ws.Range("S15").Value = 0.1830573019988039
$ws.Range("T15").Value = 0.1830573019988038
$ws.Range("E16").Value = 3
$ws.Range("G16").Value = 26.25121566666667
$ws.Range("H16").Value = 78.753647
$ws.Range("I16").Value = 0.3495454860807646
$ws.Range("J16").Value = 0.3495454860807646
$ws.Range("K16").Value = 3
$ws.Range("M16").Value = 111.1005463333333
$ws.Range("N16").Value = 333.301639
$ws.Range("O16").Value = 0.3209442197221123
$ws.Range("P16").Value = 0.3209442197221123
$ws.Range("Q16").Value = 2916.524402480826
$ws.Range("R16").Value = 26248.71962232743
$ws.Range("S16").Value = 0.1121846032875775
$ws.Range("T16").Value = 0.1121846032875775
$ws.Range("E17").Value = 3
$ws.Range("G17").Value = 26.25121566666667
$ws.Range("H17").Value = 78.753647
$ws.Range("I17").Value = 0.3495454860807646
$ws.Range("J17").Value = 0.3495454860807646
$ws.Range("K17").Value = 3
$ws.Range("M17").Value = 30.14303933333333
$ws.Range("N17").Value = 90.42911799999999
$ws.Range("O17").Value = 0.08707638763417187
$ws.Range("P17").Value = 0.08707638763417187
$ws.Range("Q17").Value = 791.2914263881495
$ws.Range("R17").Value = 7121.622837493345
$ws.Range("S17").Value = 0.03043715824174369
$ws.Range("T17").Value = 0.03043715824174369
